$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.8482917113326509
$ws.Range("F2").Value = 0.9760985490018222

$ws.Range("C3").Value = 1.025088381882265
$ws.Range("F3").Value = 1.779448698229839

$ws.Range("C4").Value = 1.743388436237755
$ws.Range("F4").Value = 2.504917496754447

$ws.Range("C5").Value = 18
$ws.Range("F5").Value = 21

$ws.Range("C6").Value = 2.030510065050192

$ws.Range("C7").Value = 1.21838349959691
$ws.Range("F7").Value = 2.056150645566329

$ws.Range("F8").Value = 0.8

$ws.Range("C9").Value = 0.4863850730369617
$ws.Range("F9").Value = 1.290076849426911
